$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values look like plain numbers need an explicit
# Text number format first, otherwise Excel auto-converts the assigned
# string into a numeric value (losing e.g. trailing zeros).
$textCells = @('D5', 'D6', 'D8', 'D12', 'D14', 'D15', 'D19', 'D21', 'D23', 'D26', 'D27', 'D28', 'D29', 'D30', 'D33', 'D34', 'D39', 'D45', 'D46', 'D48', 'D49', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '37.430.54'
$ws.Range('E2').Value = '  +0.54%  '
$ws.Range('D3').Value = '2.013.96'
$ws.Range('E3').Value = '  -0.55%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '260.06'
$ws.Range('E5').Value = '  +4.98%  '
$ws.Range('D6').Value = '0.620'
$ws.Range('E6').Value = '  -1.38%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '56.91'
$ws.Range('E8').Value = '  -5.40%  '
$ws.Range('E9').Value = '  -1.83%  '
$ws.Range('E10').Value = '  -4.51%  '
$ws.Range('E11').Value = '  -2.81%  '
$ws.Range('D12').Value = '14.31'
$ws.Range('E12').Value = '  -6.05%  '
$ws.Range('D13').Value = '2.309.61'
$ws.Range('E13').Value = '  -0.44%  '
$ws.Range('D14').Value = '21.20'
$ws.Range('E14').Value = '  -5.96%  '
$ws.Range('D15').Value = '0.803'
$ws.Range('E15').Value = '  -6.51%  '
$ws.Range('E16').Value = '  -4.87%  '
$ws.Range('D17').Value = '2.009.80'
$ws.Range('E17').Value = '  -0.70%  '
$ws.Range('D18').Value = '37.301.40'
$ws.Range('E18').Value = '  +0.46%  '
$ws.Range('D19').Value = '70.12'
$ws.Range('E19').Value = '  -0.79%  '
$ws.Range('D20').Value = '0.0₃0839'
$ws.Range('E20').Value = '  -3.35%  '
$ws.Range('D21').Value = '232.71'
$ws.Range('E21').Value = '  +0.69%  '
$ws.Range('E22').Value = '  -2.07%  '
$ws.Range('D23').Value = '2.63'
$ws.Range('E23').Value = '  +4.37%  '
$ws.Range('E24').Value = '  +0.13%  '
$ws.Range('E25').Value = '  -0.60%  '
$ws.Range('D26').Value = '164.86'
$ws.Range('E26').Value = '  +0.71%  '
$ws.Range('D27').Value = '8.95'
$ws.Range('E27').Value = '  -5.71%  '
$ws.Range('D28').Value = '19.66'
$ws.Range('E28').Value = '  -0.85%  '
$ws.Range('D29').Value = '0.131'
$ws.Range('E29').Value = '  -5.56%  '
$ws.Range('D30').Value = '1.34'
$ws.Range('E30').Value = '  -3.32%  '
$ws.Range('E31').Value = '  -1.37%  '
$ws.Range('E32').Value = '  -2.84%  '
$ws.Range('D33').Value = '4.60'
$ws.Range('E33').Value = '  -4.90%  '
$ws.Range('D34').Value = '4.55'
$ws.Range('E34').Value = '  +0.19%  '
$ws.Range('E35').Value = '  -4.13%  '
$ws.Range('E36').Value = '  +0.62%  '
$ws.Range('E37').Value = '  +0.04%  '
$ws.Range('E38').Value = '  -4.11%  '
$ws.Range('D39').Value = '5.38'
$ws.Range('E39').Value = '  -1.14%  '
$ws.Range('E40').Value = '  +3.93%  '
$ws.Range('E41').Value = '  +0.47%  '
$ws.Range('E42').Value = '  -0.95%  '
$ws.Range('E43').Value = '  -5.64%  '
$ws.Range('D44').Value = '1.414.78'
$ws.Range('E44').Value = '  +1.55%  '
$ws.Range('D45').Value = '15.84'
$ws.Range('E45').Value = '  -5.78%  '
$ws.Range('D46').Value = '90.03'
$ws.Range('E46').Value = '  -2.38%  '
$ws.Range('E47').Value = '  -3.12%  '
$ws.Range('B48').Value = 'MXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D48').Value = '2.93'
$ws.Range('E48').Value = '  +2.67%  '
$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').Value = '7.05'
$ws.Range('E49').Value = '  -5.87%  '
$ws.Range('D50').Value = '2.200.87'
$ws.Range('E50').Value = '  -0.43%  '
$ws.Range('D51').Value = '1.96'
$ws.Range('E51').Value = '  -10.54%  '
